# Update the "dSF" column (F) values for specific rows per the repull/mean
# recalculation described in the commit message.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = 1
    5  = -1
    11 = 2
    14 = 0
    16 = -1
    20 = 0
    22 = 1
    25 = 0
    27 = 5
    36 = -2
    39 = -3
    44 = 0
    51 = 0
    57 = 2
    58 = 4
    59 = -2
    64 = -1
    65 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
